$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_summaries")

$ws.Range("N2").Value = 0.7314590131410988
$ws.Range("O2").Value = 0.05826519585223427
$ws.Range("P2").Value = 0.5052730792014741
$ws.Range("Q2").Value = 0.8614241160812022
$ws.Range("R2").Value = 2.029700249935169
$ws.Range("S2").Value = 0.2453985253125837
$ws.Range("T2").Value = 1.792138120293617
$ws.Range("U2").Value = -6.505945283620632
